# Weekly update: insert a new price record as row 517 for
# "Terminal La Palmera de La Serena - Zanahoria", shifting the existing
# rows 517:536 down to 518:537.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 517; this pushes old rows 517:536 down to 518:537
$ws.Rows("517:517").Insert()

# Fill in the constant columns (same values used throughout this dataset)
$ws.Cells.Item(517, 1).Value = 8
$ws.Cells.Item(517, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(517, 3).Value = "Coquimbo"
$ws.Cells.Item(517, 4).Value = 45075
$ws.Cells.Item(517, 5).Value = 4
$ws.Cells.Item(517, 6).Value = 100114013
$ws.Cells.Item(517, 7).Value = "Zanahoria"
$ws.Cells.Item(517, 8).Value = "Sin especificar"
$ws.Cells.Item(517, 9).Value = "Primera"
$ws.Cells.Item(517, 10).Value = 500
$ws.Cells.Item(517, 11).Value = 5000
$ws.Cells.Item(517, 12).Value = 6000
$ws.Cells.Item(517, 13).Value = 5500
$ws.Cells.Item(517, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(517, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(517, 16).Value = 275
$ws.Cells.Item(517, 17).Value = 20
$ws.Cells.Item(517, 18).Value = "Hortaliza"
